$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8736128211021423
$ws.Range("B1").Value = 1.230874300003052
$ws.Range("C1").Value = 2.142632246017456
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.760517358779907
